$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Scann"
$ws.Range("B1").Value = "Scooter"

# Update data rows 2 and 3
$ws.Range("A2").Value = "https'//scooters.taxify.eu/qr/449-616"
$ws.Range("B2").Value = "449-616"

$ws.Range("A3").Value = "https'//scooters.taxify.eu/qr/449-616"
$ws.Range("B3").Value = "449-616"

# Remove now-unused rows 4 through 17
$ws.Range("A4:B17").EntireRow.Delete()
